$d = $word.ActiveDocument

# 1. Replace the lead-in text of the "Green tick image" paragraph.
$d.Content.Find.Execute("Green tick image available at: ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Completed project image, available at: ", 2)

# 2. Update the hyperlink (display text + address) that pointed at the old icon.
$h = $d.Hyperlinks.Item(3)
$h.TextToDisplay = "http://hijauangroup.com/images/completed.png"
$h.Address = "http://hijauangroup.com/images/completed.png"

# 3. Move the singleton "_GoBack" bookmark to the end of that paragraph (after
#    its trailing run, right before the paragraph mark) -- this is where Word
#    leaves it after the most recent edit.
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "Completed project image, available at:*") {
        $target = $p
        break
    }
}
$r = $target.Range
$insPos = $r.End - 1

# Bookmarks.Add chokes on a truly collapsed (zero-length) range in this
# runtime, so insert a throwaway character, wrap it with the bookmark, then
# delete the character -- the bookmark collapses in place, exactly where we
# need it.
$insRange = $d.Range($insPos, $insPos)
$insRange.InsertAfter("X")
$markRange = $d.Range($insPos, $insPos + 1)
$d.Bookmarks.Add("_GoBack", $markRange)
$markRange2 = $d.Range($insPos, $insPos + 1)
$markRange2.Delete()
